# "Generate Report for Handback" - refresh the handoff/handback timestamps
# for the en-US -> zh-cn / de-de rows that were just re-handed-back.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-26 08:53:01"
$zhcn.Range("K2").Value = "2016-08-26 08:53:27"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-26 08:53:14"
$dede.Range("K2").Value = "2016-08-26 08:53:34"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-26 08:53:14"
